$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 3976.6428  # ALC!H19: 5197.5 -> 3976.6428
$ws.Cells.Item(19, 9).Value = 5167.3  # ALC!I19: 6996.4287 -> 5167.3
$ws.Cells.Item(19, 11).Value = 5167.3  # ALC!K19: 6996.4287 -> 5167.3
$ws.Cells.Item(19, 13).Value = -4992.3  # ALC!M19: -6821.4287 -> -4992.3
$ws.Cells.Item(132, 8).Value = 1174.1637  # ALC!H132: 1114.8474 -> 1174.1637
$ws.Cells.Item(132, 9).Value = 1206.18  # ALC!I132: 1139 -> 1206.18
$ws.Cells.Item(132, 11).Value = 3618.54  # ALC!K132: 3417 -> 3618.54
$ws.Cells.Item(132, 13).Value = -1088.54  # ALC!M132: -887 -> -1088.54
$ws.Cells.Item(141, 8).Value = 18133.404  # ALC!H141: 18454.982 -> 18133.404
$ws.Cells.Item(141, 9).Value = 18165.545  # ALC!I141: 18499.629 -> 18165.545
$ws.Cells.Item(141, 11).Value = 54496.63499999999  # ALC!K141: 55498.887 -> 54496.63499999999
$ws.Cells.Item(141, 13).Value = -49316.63499999999  # ALC!M141: -50318.887 -> -49316.63499999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1890.1904  # ARM!H2: 1788.7826 -> 1890.1904
$ws.Cells.Item(2, 9).Value = 1760.2  # ARM!I2: 1625 -> 1760.2
$ws.Cells.Item(2, 10).Value = 2215.1667  # ARM!J2: 2378.4 -> 2215.1667
$ws.Cells.Item(2, 11).Value = 1760.2  # ARM!K2: 1625 -> 1760.2
$ws.Cells.Item(2, 12).Value = 2215.1667  # ARM!L2: 2378.4 -> 2215.1667
$ws.Cells.Item(2, 13).Value = -1647.2  # ARM!M2: -1512 -> -1647.2
$ws.Cells.Item(2, 14).Value = -2441.1667  # ARM!N2: -2604.4 -> -2441.1667
$ws.Cells.Item(5, 8).Value = 98.75  # ARM!H5: 84 -> 98.75
$ws.Cells.Item(5, 9).Value = 98.75  # ARM!I5: 84 -> 98.75
$ws.Cells.Item(5, 11).Value = 98.75  # ARM!K5: 84 -> 98.75
$ws.Cells.Item(5, 13).Value = 13.25  # ARM!M5: 28 -> 13.25
$ws.Cells.Item(21, 8).Value = 2591.182  # ARM!H21: 2848.2222 -> 2591.182
$ws.Cells.Item(21, 9).Value = 1085.6666  # ARM!I21: 1324 -> 1085.6666
$ws.Cells.Item(21, 10).Value = 3155.75  # ARM!J21: 3610.3333 -> 3155.75
$ws.Cells.Item(21, 11).Value = 1085.6666  # ARM!K21: 1324 -> 1085.6666
$ws.Cells.Item(21, 12).Value = 3155.75  # ARM!L21: 3610.3333 -> 3155.75
$ws.Cells.Item(21, 13).Value = -711.6666  # ARM!M21: -950 -> -711.6666
$ws.Cells.Item(21, 14).Value = -3903.75  # ARM!N21: -4358.3333 -> -3903.75
$ws.Cells.Item(74, 8).Value = 1985.75  # ARM!H74: 1856.7778 -> 1985.75
$ws.Cells.Item(74, 9).Value = 1878.5333  # ARM!I74: 1791.625 -> 1878.5333
$ws.Cells.Item(74, 10).Value = 2164.4443  # ARM!J74: 1951.5454 -> 2164.4443
$ws.Cells.Item(74, 11).Value = 1878.5333  # ARM!K74: 1791.625 -> 1878.5333
$ws.Cells.Item(74, 12).Value = 2164.4443  # ARM!L74: 1951.5454 -> 2164.4443
$ws.Cells.Item(74, 13).Value = -1004.5333  # ARM!M74: -917.625 -> -1004.5333
$ws.Cells.Item(74, 14).Value = -3912.4443  # ARM!N74: -3699.5454 -> -3912.4443
$ws.Cells.Item(77, 8).Value = 1985.75  # ARM!H77: 1856.7778 -> 1985.75
$ws.Cells.Item(77, 9).Value = 1878.5333  # ARM!I77: 1791.625 -> 1878.5333
$ws.Cells.Item(77, 10).Value = 2164.4443  # ARM!J77: 1951.5454 -> 2164.4443
$ws.Cells.Item(77, 11).Value = 9392.666499999999  # ARM!K77: 8958.125 -> 9392.666499999999
$ws.Cells.Item(77, 12).Value = 10822.2215  # ARM!L77: 9757.726999999999 -> 10822.2215
$ws.Cells.Item(77, 13).Value = -5024.666499999999  # ARM!M77: -4590.125 -> -5024.666499999999
$ws.Cells.Item(77, 14).Value = -19558.2215  # ARM!N77: -18493.727 -> -19558.2215
$ws.Cells.Item(116, 8).Value = 1890.1904  # ARM!H116: 1788.7826 -> 1890.1904
$ws.Cells.Item(116, 9).Value = 1760.2  # ARM!I116: 1625 -> 1760.2
$ws.Cells.Item(116, 10).Value = 2215.1667  # ARM!J116: 2378.4 -> 2215.1667
$ws.Cells.Item(116, 11).Value = 1760.2  # ARM!K116: 1625 -> 1760.2
$ws.Cells.Item(116, 12).Value = 2215.1667  # ARM!L116: 2378.4 -> 2215.1667
$ws.Cells.Item(116, 13).Value = 533.8  # ARM!M116: 669 -> 533.8
$ws.Cells.Item(116, 14).Value = -6803.1667  # ARM!N116: -6966.4 -> -6803.1667
$ws.Cells.Item(122, 8).Value = 1587.0465  # ARM!H122: 1680.1025 -> 1587.0465
$ws.Cells.Item(122, 9).Value = 1529.3948  # ARM!I122: 1610.2059 -> 1529.3948
$ws.Cells.Item(122, 10).Value = 2025.2  # ARM!J122: 2155.4 -> 2025.2
$ws.Cells.Item(122, 11).Value = 4588.1844  # ARM!K122: 4830.6177 -> 4588.1844
$ws.Cells.Item(122, 12).Value = 6075.6  # ARM!L122: 6466.200000000001 -> 6075.6
$ws.Cells.Item(122, 13).Value = -2138.1844  # ARM!M122: -2380.6177 -> -2138.1844
$ws.Cells.Item(122, 14).Value = -10975.6  # ARM!N122: -11366.2 -> -10975.6
$ws.Cells.Item(132, 8).Value = 1930.875  # ARM!H132: 1301.5957 -> 1930.875
$ws.Cells.Item(132, 9).Value = 1791.5  # ARM!I132: 1205.4889 -> 1791.5
$ws.Cells.Item(132, 11).Value = 5374.5  # ARM!K132: 3616.4667 -> 5374.5
$ws.Cells.Item(132, 13).Value = -2844.5  # ARM!M132: -1086.4667 -> -2844.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1890.1904  # BSM!H3: 1788.7826 -> 1890.1904
$ws.Cells.Item(3, 9).Value = 1760.2  # BSM!I3: 1625 -> 1760.2
$ws.Cells.Item(3, 10).Value = 2215.1667  # BSM!J3: 2378.4 -> 2215.1667
$ws.Cells.Item(3, 11).Value = 1760.2  # BSM!K3: 1625 -> 1760.2
$ws.Cells.Item(3, 12).Value = 2215.1667  # BSM!L3: 2378.4 -> 2215.1667
$ws.Cells.Item(3, 13).Value = -1646.2  # BSM!M3: -1511 -> -1646.2
$ws.Cells.Item(3, 14).Value = -2443.1667  # BSM!N3: -2606.4 -> -2443.1667
$ws.Cells.Item(4, 8).Value = 98.75  # BSM!H4: 84 -> 98.75
$ws.Cells.Item(4, 9).Value = 98.75  # BSM!I4: 84 -> 98.75
$ws.Cells.Item(4, 11).Value = 98.75  # BSM!K4: 84 -> 98.75
$ws.Cells.Item(4, 13).Value = 16.25  # BSM!M4: 31 -> 16.25
$ws.Cells.Item(5, 8).Value = 1000  # BSM!H5: 401.33334 -> 1000
$ws.Cells.Item(5, 9).Value = 0  # BSM!I5: 102 -> 0
$ws.Cells.Item(5, 11).Value = 0  # BSM!K5: 102 -> 0
$ws.Cells.Item(5, 13).ClearContents()  # BSM!M5 was 11
$ws.Cells.Item(80, 8).Value = 1058.8334  # BSM!H80: 1027.9412 -> 1058.8334
$ws.Cells.Item(80, 9).Value = 729.1667  # BSM!I80: 638.5714 -> 729.1667
$ws.Cells.Item(80, 10).Value = 1223.6666  # BSM!J80: 1300.5 -> 1223.6666
$ws.Cells.Item(80, 11).Value = 729.1667  # BSM!K80: 638.5714 -> 729.1667
$ws.Cells.Item(80, 12).Value = 1223.6666  # BSM!L80: 1300.5 -> 1223.6666
$ws.Cells.Item(80, 13).Value = 268.8333  # BSM!M80: 359.4286 -> 268.8333
$ws.Cells.Item(80, 14).Value = -3219.6666  # BSM!N80: -3296.5 -> -3219.6666
$ws.Cells.Item(83, 8).Value = 1058.8334  # BSM!H83: 1027.9412 -> 1058.8334
$ws.Cells.Item(83, 9).Value = 729.1667  # BSM!I83: 638.5714 -> 729.1667
$ws.Cells.Item(83, 10).Value = 1223.6666  # BSM!J83: 1300.5 -> 1223.6666
$ws.Cells.Item(83, 11).Value = 3645.8335  # BSM!K83: 3192.857 -> 3645.8335
$ws.Cells.Item(83, 12).Value = 6118.333000000001  # BSM!L83: 6502.5 -> 6118.333000000001
$ws.Cells.Item(83, 13).Value = 1346.1665  # BSM!M83: 1799.143 -> 1346.1665
$ws.Cells.Item(83, 14).Value = -16102.333  # BSM!N83: -16486.5 -> -16102.333
$ws.Cells.Item(107, 8).Value = 20638.434  # BSM!H107: 21032.674 -> 20638.434
$ws.Cells.Item(107, 9).Value = 26969.46  # BSM!I107: 27675.553 -> 26969.46
$ws.Cells.Item(107, 11).Value = 26969.46  # BSM!K107: 27675.553 -> 26969.46
$ws.Cells.Item(107, 13).Value = -25049.46  # BSM!M107: -25755.553 -> -25049.46
$ws.Cells.Item(134, 8).Value = 905.5238000000001  # BSM!H134: 941.8 -> 905.5238000000001
$ws.Cells.Item(134, 9).Value = 905.5238000000001  # BSM!I134: 941.8 -> 905.5238000000001
$ws.Cells.Item(134, 11).Value = 2716.5714  # BSM!K134: 2825.4 -> 2716.5714
$ws.Cells.Item(134, 13).Value = -181.5714000000003  # BSM!M134: -290.3999999999996 -> -181.5714000000003
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3854.5715  # CRP!H16: 3497.75 -> 3854.5715
$ws.Cells.Item(16, 9).Value = 4196.6  # CRP!I16: 3663.8333 -> 4196.6
$ws.Cells.Item(16, 11).Value = 4196.6  # CRP!K16: 3663.8333 -> 4196.6
$ws.Cells.Item(16, 13).Value = -3909.6  # CRP!M16: -3376.8333 -> -3909.6
$ws.Cells.Item(31, 8).Value = 48293.4  # CRP!H31: 49824.207 -> 48293.4
$ws.Cells.Item(31, 9).Value = 104379.9  # CRP!I31: 115544.336 -> 104379.9
$ws.Cells.Item(31, 11).Value = 104379.9  # CRP!K31: 115544.336 -> 104379.9
$ws.Cells.Item(31, 13).Value = -104084.9  # CRP!M31: -115249.336 -> -104084.9
$ws.Cells.Item(34, 8).Value = 48293.4  # CRP!H34: 49824.207 -> 48293.4
$ws.Cells.Item(34, 9).Value = 104379.9  # CRP!I34: 115544.336 -> 104379.9
$ws.Cells.Item(34, 11).Value = 104379.9  # CRP!K34: 115544.336 -> 104379.9
$ws.Cells.Item(34, 13).Value = -104177.9  # CRP!M34: -115342.336 -> -104177.9
$ws.Cells.Item(43, 8).Value = 17551.166  # CRP!H43: 15615.286 -> 17551.166
$ws.Cells.Item(43, 10).Value = 17551.166  # CRP!J43: 15615.286 -> 17551.166
$ws.Cells.Item(43, 12).Value = 17551.166  # CRP!L43: 15615.286 -> 17551.166
$ws.Cells.Item(43, 14).Value = -17919.166  # CRP!N43: -15983.286 -> -17919.166
$ws.Cells.Item(59, 8).Value = 19999.777  # CRP!H59: 19999.875 -> 19999.777
$ws.Cells.Item(59, 10).Value = 19999.777  # CRP!J59: 19999.875 -> 19999.777
$ws.Cells.Item(59, 12).Value = 19999.777  # CRP!L59: 19999.875 -> 19999.777
$ws.Cells.Item(59, 14).Value = -22289.777  # CRP!N59: -22289.875 -> -22289.777
$ws.Cells.Item(101, 8).Value = 17551.166  # CRP!H101: 15615.286 -> 17551.166
$ws.Cells.Item(101, 10).Value = 17551.166  # CRP!J101: 15615.286 -> 17551.166
$ws.Cells.Item(101, 12).Value = 17551.166  # CRP!L101: 15615.286 -> 17551.166
$ws.Cells.Item(101, 14).Value = -24041.166  # CRP!N101: -22105.286 -> -24041.166
$ws.Cells.Item(113, 8).Value = 3854.5715  # CRP!H113: 3497.75 -> 3854.5715
$ws.Cells.Item(113, 9).Value = 4196.6  # CRP!I113: 3663.8333 -> 4196.6
$ws.Cells.Item(113, 11).Value = 4196.6  # CRP!K113: 3663.8333 -> 4196.6
$ws.Cells.Item(113, 13).Value = -2026.6  # CRP!M113: -1493.8333 -> -2026.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 2500369.2  # CUL!H50: 2500094.2 -> 2500369.2
$ws.Cells.Item(50, 9).Value = 443.2  # CUL!I50: 113.2 -> 443.2
$ws.Cells.Item(50, 11).Value = 1329.6  # CUL!K50: 339.6 -> 1329.6
$ws.Cells.Item(50, 13).Value = -848.5999999999999  # CUL!M50: 141.4 -> -848.5999999999999
$ws.Cells.Item(53, 8).Value = 2500369.2  # CUL!H53: 2500094.2 -> 2500369.2
$ws.Cells.Item(53, 9).Value = 443.2  # CUL!I53: 113.2 -> 443.2
$ws.Cells.Item(53, 11).Value = 1329.6  # CUL!K53: 339.6 -> 1329.6
$ws.Cells.Item(53, 13).Value = -848.5999999999999  # CUL!M53: 141.4 -> -848.5999999999999
$ws.Cells.Item(58, 8).Value = 9215  # CUL!H58: 5254 -> 9215
$ws.Cells.Item(58, 9).Value = 5391.6665  # CUL!I58: 4087.5 -> 5391.6665
$ws.Cells.Item(58, 10).Value = 14950  # CUL!J58: 7587 -> 14950
$ws.Cells.Item(58, 11).Value = 16174.9995  # CUL!K58: 12262.5 -> 16174.9995
$ws.Cells.Item(58, 12).Value = 44850  # CUL!L58: 22761 -> 44850
$ws.Cells.Item(58, 13).Value = -16046.9995  # CUL!M58: -12134.5 -> -16046.9995
$ws.Cells.Item(58, 14).Value = -45106  # CUL!N58: -23017 -> -45106
$ws.Cells.Item(131, 8).Value = 1567.1875  # CUL!H131: 1578.3572 -> 1567.1875
$ws.Cells.Item(131, 9).Value = 1024.4  # CUL!I131: 1053.6666 -> 1024.4
$ws.Cells.Item(131, 10).Value = 1813.909  # CUL!J131: 1971.875 -> 1813.909
$ws.Cells.Item(131, 11).Value = 3073.2  # CUL!K131: 3160.9998 -> 3073.2
$ws.Cells.Item(131, 12).Value = 5441.727000000001  # CUL!L131: 5915.625 -> 5441.727000000001
$ws.Cells.Item(131, 13).Value = 1966.8  # CUL!M131: 1879.0002 -> 1966.8
$ws.Cells.Item(131, 14).Value = -15521.727  # CUL!N131: -15995.625 -> -15521.727
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5666.6665  # GSM!H70: 5699 -> 5666.6665
$ws.Cells.Item(70, 9).Value = 5600  # GSM!I70: 5682.1665 -> 5600
$ws.Cells.Item(70, 11).Value = 5600  # GSM!K70: 5682.1665 -> 5600
$ws.Cells.Item(70, 13).Value = -5330  # GSM!M70: -5412.1665 -> -5330
$ws.Cells.Item(73, 8).Value = 5666.6665  # GSM!H73: 5699 -> 5666.6665
$ws.Cells.Item(73, 9).Value = 5600  # GSM!I73: 5682.1665 -> 5600
$ws.Cells.Item(73, 11).Value = 5600  # GSM!K73: 5682.1665 -> 5600
$ws.Cells.Item(73, 13).Value = -4664  # GSM!M73: -4746.1665 -> -4664
$ws.Cells.Item(97, 8).Value = 18944.578  # GSM!H97: 20028.055 -> 18944.578
$ws.Cells.Item(97, 9).Value = 24110.518  # GSM!I97: 24179.965 -> 24110.518
$ws.Cells.Item(97, 10).Value = 2298.7778  # GSM!J97: 2827.2856 -> 2298.7778
$ws.Cells.Item(97, 11).Value = 24110.518  # GSM!K97: 24179.965 -> 24110.518
$ws.Cells.Item(97, 12).Value = 2298.7778  # GSM!L97: 2827.2856 -> 2298.7778
$ws.Cells.Item(97, 13).Value = -23614.518  # GSM!M97: -23683.965 -> -23614.518
$ws.Cells.Item(97, 14).Value = -3290.7778  # GSM!N97: -3819.2856 -> -3290.7778
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1007.8571  # LTW!H22: 796.6667 -> 1007.8571
$ws.Cells.Item(22, 9).Value = 1173.5  # LTW!I22: 822.2143 -> 1173.5
$ws.Cells.Item(22, 10).Value = 787  # LTW!J22: 745.5714 -> 787
$ws.Cells.Item(22, 11).Value = 1173.5  # LTW!K22: 822.2143 -> 1173.5
$ws.Cells.Item(22, 12).Value = 787  # LTW!L22: 745.5714 -> 787
$ws.Cells.Item(22, 13).Value = -878.5  # LTW!M22: -527.2143 -> -878.5
$ws.Cells.Item(22, 14).Value = -1377  # LTW!N22: -1335.5714 -> -1377
$ws.Cells.Item(27, 8).Value = 1007.8571  # LTW!H27: 796.6667 -> 1007.8571
$ws.Cells.Item(27, 9).Value = 1173.5  # LTW!I27: 822.2143 -> 1173.5
$ws.Cells.Item(27, 10).Value = 787  # LTW!J27: 745.5714 -> 787
$ws.Cells.Item(27, 11).Value = 1173.5  # LTW!K27: 822.2143 -> 1173.5
$ws.Cells.Item(27, 12).Value = 787  # LTW!L27: 745.5714 -> 787
$ws.Cells.Item(27, 13).Value = -1066.5  # LTW!M27: -715.2143 -> -1066.5
$ws.Cells.Item(27, 14).Value = -1001  # LTW!N27: -959.5714 -> -1001
$ws.Cells.Item(46, 8).Value = 25044.389  # LTW!H46: 18977.084 -> 25044.389
$ws.Cells.Item(46, 9).Value = 212260  # LTW!I46: 85374 -> 212260
$ws.Cells.Item(46, 10).Value = 1642.4375  # LTW!J46: 1504.2106 -> 1642.4375
$ws.Cells.Item(46, 11).Value = 212260  # LTW!K46: 85374 -> 212260
$ws.Cells.Item(46, 12).Value = 1642.4375  # LTW!L46: 1504.2106 -> 1642.4375
$ws.Cells.Item(46, 13).Value = -212072  # LTW!M46: -85186 -> -212072
$ws.Cells.Item(46, 14).Value = -2018.4375  # LTW!N46: -1880.2106 -> -2018.4375
$ws.Cells.Item(55, 8).Value = 350.10526  # LTW!H55: 362.94446 -> 350.10526
$ws.Cells.Item(55, 10).Value = 403.77777  # LTW!J55: 439.375 -> 403.77777
$ws.Cells.Item(55, 12).Value = 403.77777  # LTW!L55: 439.375 -> 403.77777
$ws.Cells.Item(55, 14).Value = -749.7777699999999  # LTW!N55: -785.375 -> -749.7777699999999
$ws.Cells.Item(100, 8).Value = 19538.6  # LTW!H100: 17989.363 -> 19538.6
$ws.Cells.Item(100, 9).Value = 1967.3636  # LTW!I100: 2024.2 -> 1967.3636
$ws.Cells.Item(100, 10).Value = 41014.555  # LTW!J100: 31293.666 -> 41014.555
$ws.Cells.Item(100, 11).Value = 1967.3636  # LTW!K100: 2024.2 -> 1967.3636
$ws.Cells.Item(100, 12).Value = 41014.555  # LTW!L100: 31293.666 -> 41014.555
$ws.Cells.Item(100, 13).Value = -1426.3636  # LTW!M100: -1483.2 -> -1426.3636
$ws.Cells.Item(100, 14).Value = -42096.555  # LTW!N100: -32375.666 -> -42096.555
$ws.Cells.Item(136, 8).Value = 4879.8  # LTW!H136: 4774.8125 -> 4879.8
$ws.Cells.Item(136, 9).Value = 4321.2  # LTW!I136: 4219.273 -> 4321.2
$ws.Cells.Item(136, 11).Value = 12963.6  # LTW!K136: 12657.819 -> 12963.6
$ws.Cells.Item(136, 13).Value = -10413.6  # LTW!M136: -10107.819 -> -10413.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(115, 8).Value = 80000  # WVR!H115: 80377 -> 80000
$ws.Cells.Item(115, 10).Value = 80000  # WVR!J115: 80377 -> 80000
$ws.Cells.Item(115, 12).Value = 80000  # WVR!L115: 80377 -> 80000
$ws.Cells.Item(115, 14).Value = -83134  # WVR!N115: -83511 -> -83134
$ws.Cells.Item(135, 8).Value = 89257.5  # WVR!H135: 82504.664 -> 89257.5
$ws.Cells.Item(135, 10).Value = 89257.5  # WVR!J135: 82504.664 -> 89257.5
$ws.Cells.Item(135, 12).Value = 89257.5  # WVR!L135: 82504.664 -> 89257.5
$ws.Cells.Item(135, 14).Value = -99397.5  # WVR!N135: -92644.664 -> -99397.5
$ws.Cells.Item(141, 8).Value = 75994  # WVR!H141: 75993.75 -> 75994
$ws.Cells.Item(141, 10).Value = 75994  # WVR!J141: 75993.75 -> 75994
$ws.Cells.Item(141, 12).Value = 75994  # WVR!L141: 75993.75 -> 75994
$ws.Cells.Item(141, 14).Value = -86354  # WVR!N141: -86353.75 -> -86354
